$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new, most-recent week's record is inserted at row 6,
# pushing the previous rows 6-9 down to rows 7-10 (dimension grows to A1:R10).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44580
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100114007
$ws.Range("G6").Value = "Jengibre"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11500
$ws.Range("N6").Value = '$/caja 13 kilos'
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 885
$ws.Range("Q6").Value = 13
$ws.Range("R6").Value = "Hortaliza"
